$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "Due Date" column (column J) is being removed entirely.
# Deleting the whole column shifts everything after it one column to the left.
$ws.Columns.Item(10).Delete()

# Update the selection to match what a user would see after selecting and deleting column J
$ws.Range("J1:J1048576").Select()
